$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 359 (shifts existing rows 359..435 down to 360..436)
$ws.Rows.Item(359).Insert()

# Fill in the new row 359 with the new data point
$ws.Cells.Item(359, 1).Value = 9
$ws.Cells.Item(359, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(359, 3).Value = "Metropolitana"
$ws.Cells.Item(359, 4).Value = 44711
$ws.Cells.Item(359, 4).NumberFormat = $ws.Cells.Item(360, 4).NumberFormat
$ws.Cells.Item(359, 5).Value = 13
$ws.Cells.Item(359, 6).Value = 100112012
$ws.Cells.Item(359, 7).Value = "Espinaca"
$ws.Cells.Item(359, 8).Value = "Sin especificar"
$ws.Cells.Item(359, 9).Value = "Primera"
$ws.Cells.Item(359, 10).Value = 70
$ws.Cells.Item(359, 11).Value = 6000
$ws.Cells.Item(359, 12).Value = 6000
$ws.Cells.Item(359, 13).Value = 6000
$ws.Cells.Item(359, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(359, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(359, 16).Value = 600
$ws.Cells.Item(359, 17).Value = 10
$ws.Cells.Item(359, 18).Value = "Hortaliza"
